$wb = $excel.ActiveWorkbook

# Helper: write a value into a cell as TEXT (not a number), without leaving
# a custom number format behind on the cell (matches source cells that were
# already stored as text/inline-strings elsewhere in the workbook).
function Set-TextValue {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# ---- Overall sheet ----
$wsOverall = $wb.Worksheets.Item("Overall")
Set-TextValue $wsOverall "A2" "663"

# ---- County sheet ----
$wsCounty = $wb.Worksheets.Item("County")
Set-TextValue $wsCounty "B2" "29"
Set-TextValue $wsCounty "B3" "55"
Set-TextValue $wsCounty "B4" "39"
Set-TextValue $wsCounty "B5" "18"
Set-TextValue $wsCounty "B6" "75"
Set-TextValue $wsCounty "B7" "162"
Set-TextValue $wsCounty "B8" "118"
Set-TextValue $wsCounty "B9" "101"
Set-TextValue $wsCounty "B10" "48"
Set-TextValue $wsCounty "B11" "18"

# New "Total" row for County sheet
Set-TextValue $wsCounty "A12" "Total"
Set-TextValue $wsCounty "B12" "663"
Set-TextValue $wsCounty "C12" "`$1,133,383,584"
Set-TextValue $wsCounty "D12" "10.55%"
Set-TextValue $wsCounty "E12" "-7.39%"
Set-TextValue $wsCounty "F12" "62.90%"

# ---- Congressional District sheet ----
$wsCd = $wb.Worksheets.Item("Congressional District")
Set-TextValue $wsCd "B2" "319"
Set-TextValue $wsCd "B3" "344"
Set-TextValue $wsCd "B4" "663"

# ---- Size sheet ----
$wsSize = $wb.Worksheets.Item("Size")
Set-TextValue $wsSize "B2" "227"
Set-TextValue $wsSize "B3" "176"
Set-TextValue $wsSize "B4" "113"
Set-TextValue $wsSize "B5" "27"
Set-TextValue $wsSize "B6" "95"
Set-TextValue $wsSize "B7" "25"
Set-TextValue $wsSize "B8" "663"

# ---- Subsector sheet ----
$wsSub = $wb.Worksheets.Item("Subsector")
Set-TextValue $wsSub "B2" "54"
Set-TextValue $wsSub "B3" "86"
Set-TextValue $wsSub "B4" "36"
Set-TextValue $wsSub "B5" "64"
Set-TextValue $wsSub "B6" "4"
Set-TextValue $wsSub "B7" "211"
Set-TextValue $wsSub "B8" "3"
Set-TextValue $wsSub "B9" "40"
Set-TextValue $wsSub "B10" "6"
Set-TextValue $wsSub "B11" "150"
Set-TextValue $wsSub "B12" "9"
Set-TextValue $wsSub "B13" "663"
